# Apply the edit described by the diff:
# A new data record was inserted at row 434 (pushing all subsequent
# records down by one row, through the former last row 529 which is
# now row 530). The sheet dimension grows from A1:R529 to A1:R530.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 434; Excel shifts rows 434..529 down to 435..530
# and copies formatting (incl. the date number format on column D) from the row above.
$ws.Rows.Item(434).Insert()

# Populate the newly inserted row 434 with the new record's data.
$ws.Range("A434").Value = 10
$ws.Range("B434").Value = "Vega Modelo de Temuco"
$ws.Range("C434").Value = "La Araucanía"
$ws.Range("D434").Value = 44995
$ws.Range("E434").Value = 9
$ws.Range("F434").Value = 100112040
$ws.Range("G434").Value = "Cilantro"
$ws.Range("H434").Value = "Sin especificar"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 55
$ws.Range("K434").Value = 7000
$ws.Range("L434").Value = 7000
$ws.Range("M434").Value = 7000
$ws.Range("N434").Value = "$/docena de atados (2 kilos)"
$ws.Range("O434").Value = "Provincia de Cautín"
$ws.Range("P434").Value = 3500
$ws.Range("Q434").Value = 2
$ws.Range("R434").Value = "Hortaliza"
